$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column G, matching style of existing header cells (A1:F1)
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Baseline Accuracy"

# Fill in baseline accuracy values for each row
$ws.Range("G2").Value = 0.7129947253882188
$ws.Range("G3").Value = 0.6596215558260484
$ws.Range("G4").Value = 0.6181992549149792
$ws.Range("G5").Value = 0.5794327025930435
$ws.Range("G6").Value = 0.5145144037475563
$ws.Range("G7").Value = 0.5118586551584228
$ws.Range("G8").Value = 0.6250599387702409
$ws.Range("G9").Value = 0.6377485153627679
$ws.Range("G10").Value = 0.6577772859724835
$ws.Range("G11").Value = 0.8182287632326362
$ws.Range("G12").Value = 0.8931430046844454
